$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the review completion date in F6 (DATE column of the new preparation row)
$ws.Range("F6").Value = 45224

# Item #17 (row 28) moved from "en cours" / "Dina" to "fait" / "Elie" -> apply the
# same "done" (green) formatting used by the other completed rows (e.g. row 27/29)
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("F27").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G27").Copy()
$ws.Range("G28").PasteSpecial(-4122)

$ws.Range("F28").Value = "Elie"
$ws.Range("G28").Value = "fait"

# Move the active selection to I13, matching the reviewer's new cursor position
$ws.Range("I13").Select()
